$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Row 1 headers A:F, then row 2 data A:F (matches original authoring order) ---
$ws2.Range("A1").Value = "Profile_FirstName"
$ws2.Range("B1").Value = "profile_LastName"
$ws2.Range("C1").Value = "profile_Email"
$ws2.Range("D1").Value = "profile_Password"
$ws2.Range("E1").Value = "profile_ConfirmPassword"
$ws2.Range("F1").Value = "profile_ZipCode"

$ws2.Range("A2").Value = "venkat"
$ws2.Range("B2").Value = "mareedu"

$ws2.Range("C2").Value = "venkat@cgi.com"
$ws2.Hyperlinks.Add($ws2.Range("C2"), "mailto:venkat@cgi.com")
$ws2.Range("C2").Style = "Hyperlink"

$ws2.Range("D2").Value = "dada@123"
$ws2.Hyperlinks.Add($ws2.Range("D2"), "mailto:dada@123")
$ws2.Range("D2").Style = "Hyperlink"

$ws2.Range("E2").Value = "dada@123"
$ws2.Hyperlinks.Add($ws2.Range("E2"), "mailto:dada@123")
$ws2.Range("E2").Style = "Hyperlink"

$ws2.Range("F2").Value = 76502

# --- Remaining columns G:K, header then data per column ---
$ws2.Range("G1").Value = "profile_PhNumber"
$ws2.Range("G2").Value = "'1111111111"

$ws2.Range("H1").Value = "profile_inValidPhNumber"
$ws2.Range("H2").Value = "'111111111"

$ws2.Range("I1").Value = "profile_inValidemail"
$ws2.Range("I2").Value = "venka"

$ws2.Range("J1").Value = "profile_inValidezip"
$ws2.Range("J2").Value = "'7645"

$ws2.Range("K1").Value = "profile_InvalidPassword"
$ws2.Range("K2").Value = "dada@1234"
$ws2.Hyperlinks.Add($ws2.Range("K2"), "mailto:dada@1234")
$ws2.Range("K2").Style = "Hyperlink"

# --- Header fill color (green) across A1:K1 ---
$ws2.Range("A1:K1").Interior.Color = 5287936

# --- Column widths (best-fit, derived from the rendered content) ---
$ws2.Columns.Item(1).ColumnWidth = 14.944010416666666
$ws2.Columns.Item(2).ColumnWidth = 14.830729166666666
$ws2.Columns.Item(3).ColumnWidth = 13.830729166666666
$ws2.Columns.Item(4).ColumnWidth = 14.385416666666666
$ws2.Columns.Item(5).ColumnWidth = 21.166666666666668
$ws2.Columns.Item(6).ColumnWidth = 13.166666666666666
$ws2.Columns.Item(7).ColumnWidth = 15.166666666666666
$ws2.Columns.Item(8).ColumnWidth = 20.830729166666668
$ws2.Columns.Item(9).ColumnWidth = 20.830729166666668
$ws2.Columns.Item(10).ColumnWidth = 15.276041666666666
$ws2.Columns.Item(11).ColumnWidth = 19.830729166666668

# --- Selections: sheet1 selection moves to C1 (it is no longer the active tab) ---
$ws1.Range("C1").Select() | Out-Null

# --- Switch to sheet2 as the active tab, with final selection F21 ---
$ws2.Activate() | Out-Null
$ws2.Range("F21").Select() | Out-Null
